# Updates the Price (D) and Volume(1h) (E) columns for the cryptos
# table on row 2-51 with refreshed market data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD = $ws.Range("D2")
$cD.Value2 = "'24.000.87"
$cD.ClearFormats()
$ws.Range("E2").Value = '  +1.05%  '

$cD = $ws.Range("D3")
$cD.Value2 = "'1.654.09"
$cD.ClearFormats()
$ws.Range("E3").Value = '  +2.37%  '

$cD = $ws.Range("D4")
$cD.Value2 = "'0.9982"
$cD.ClearFormats()
$ws.Range("E4").Value = '  -0.49%  '

$cD = $ws.Range("D5")
$cD.Value2 = "'308.76"
$cD.ClearFormats()
$ws.Range("E5").Value = '  +0.69%  '

$cD = $ws.Range("D6")
$cD.Value2 = "'0.9987"
$cD.ClearFormats()
$ws.Range("E6").Value = '  -0.17%  '

$ws.Range("E7").Value = '  -0.42%  '

$cD = $ws.Range("D8")
$cD.Value2 = "'0.3832"
$cD.ClearFormats()
$ws.Range("E8").Value = '  -0.10%  '

$cD = $ws.Range("D9")
$cD.Value2 = "'51.26"
$cD.ClearFormats()
$ws.Range("E9").Value = '  +3.13%  '

$cD = $ws.Range("D10")
$cD.Value2 = "'1.357"
$cD.ClearFormats()
$ws.Range("E10").Value = '  -0.38%  '

$cD = $ws.Range("D11")
$cD.Value2 = "'0.9979"
$cD.ClearFormats()
$ws.Range("E11").Value = '  -0.57%  '

$cD = $ws.Range("D12")
$cD.Value2 = "'0.08444"
$cD.ClearFormats()
$ws.Range("E12").Value = '  +0.21%  '

$cD = $ws.Range("D13")
$cD.Value2 = "'24.01"
$cD.ClearFormats()
$ws.Range("E13").Value = '  +0.98%  '

$cD = $ws.Range("D14")
$cD.Value2 = "'7.127"
$cD.ClearFormats()
$ws.Range("E14").Value = '  +1.34%  '

$cD = $ws.Range("D15")
$cD.Value2 = "'7.881"
$cD.ClearFormats()
$ws.Range("E15").Value = '  +4.74%  '

$ws.Range("E16").Value = '  +3.26%  '

$cD = $ws.Range("D17")
$cD.Value2 = "'1.651.05"
$cD.ClearFormats()

$cD = $ws.Range("D18")
$cD.Value2 = "'94.48"
$cD.ClearFormats()
$ws.Range("E18").Value = '  +0.79%  '

$cD = $ws.Range("D19")
$cD.Value2 = "'0.06971"
$cD.ClearFormats()
$ws.Range("E19").Value = '  +0.78%  '

$cD = $ws.Range("D20")
$cD.Value2 = "'19.82"
$cD.ClearFormats()
$ws.Range("E20").Value = '  -0.74%  '

$cD = $ws.Range("D21")
$cD.Value2 = "'6.895"
$cD.ClearFormats()
$ws.Range("E21").Value = '  +1.37%  '

$cD = $ws.Range("D22")
$cD.Value2 = "'0.9985"
$cD.ClearFormats()
$ws.Range("E22").Value = '  -0.18%  '

$cD = $ws.Range("D23")
$cD.Value2 = "'13.63"
$cD.ClearFormats()
$ws.Range("E23").Value = '  +2.05%  '

$cD = $ws.Range("D24")
$cD.Value2 = "'23.987.52"
$cD.ClearFormats()
$ws.Range("E24").Value = '  +0.93%  '

$ws.Range("E25").Value = '  +1.36%  '

$cD = $ws.Range("D26")
$cD.Value2 = "'3.036"
$cD.ClearFormats()
$ws.Range("E26").Value = '  +6.51%  '

$cD = $ws.Range("D27")
$cD.Value2 = "'22.07"
$cD.ClearFormats()
$ws.Range("E27").Value = '  -0.57%  '

$cD = $ws.Range("D28")
$cD.Value2 = "'152.76"
$cD.ClearFormats()
$ws.Range("E28").Value = '  -2.21%  '

$cD = $ws.Range("D29")
$cD.Value2 = "'5.438"
$cD.ClearFormats()
$ws.Range("E29").Value = '  +3.49%  '

$cD = $ws.Range("D30")
$cD.Value2 = "'139.12"
$cD.ClearFormats()
$ws.Range("E30").Value = '  -0.56%  '

$cD = $ws.Range("D31")
$cD.Value2 = "'7.763"
$cD.ClearFormats()
$ws.Range("E31").Value = '  -0.91%  '

$cD = $ws.Range("D32")
$cD.Value2 = "'2.484"
$cD.ClearFormats()
$ws.Range("E32").Value = '  -0.56%  '

$cD = $ws.Range("D33")
$cD.Value2 = "'1.834.05"
$cD.ClearFormats()
$ws.Range("E33").Value = '  +2.06%  '

$cD = $ws.Range("D34")
$cD.Value2 = "'1.037"
$cD.ClearFormats()
$ws.Range("E34").Value = '  +6.00%  '

$cD = $ws.Range("D35")
$cD.Value2 = "'0.08131"
$cD.ClearFormats()
$ws.Range("E35").Value = '  +0.21%  '

$cD = $ws.Range("D36")
$cD.Value2 = "'0.02971"
$cD.ClearFormats()
$ws.Range("E36").Value = '  +3.38%  '

$cD = $ws.Range("D37")
$cD.Value2 = "'6.728"
$cD.ClearFormats()
$ws.Range("E37").Value = '  +2.52%  '

$cD = $ws.Range("D38")
$cD.Value2 = "'10.88"
$cD.ClearFormats()
$ws.Range("E38").Value = '  +6.00%  '

$cD = $ws.Range("D39")
$cD.Value2 = "'0.2684"
$cD.ClearFormats()
$ws.Range("E39").Value = '  +0.81%  '

$cD = $ws.Range("D40")
$cD.Value2 = "'0.09138"
$cD.ClearFormats()
$ws.Range("E40").Value = '  +0.16%  '

$cD = $ws.Range("D41")
$cD.Value2 = "'0.7567"
$cD.ClearFormats()
$ws.Range("E41").Value = '  +1.15%  '

$cD = $ws.Range("D42")
$cD.Value2 = "'13.49"
$cD.ClearFormats()
$ws.Range("E42").Value = '  +0.09%  '

$cD = $ws.Range("D43")
$cD.Value2 = "'1.429"
$cD.ClearFormats()
$ws.Range("E43").Value = '  +0.69%  '

$cD = $ws.Range("D44")
$cD.Value2 = "'16.33"
$cD.ClearFormats()
$ws.Range("E44").Value = '  +1.83%  '

$cD = $ws.Range("D45")
$cD.Value2 = "'0.6947"
$cD.ClearFormats()
$ws.Range("E45").Value = '  +1.06%  '

$cD = $ws.Range("D46")
$cD.Value2 = "'2.456"
$cD.ClearFormats()
$ws.Range("E46").Value = '  -0.20%  '

$cD = $ws.Range("D47")
$cD.Value2 = "'4.074"
$cD.ClearFormats()
$ws.Range("E47").Value = '  +0.38%  '

$cD = $ws.Range("D48")
$cD.Value2 = "'0.9979"
$cD.ClearFormats()
$ws.Range("E48").Value = '  -0.25%  '

$cD = $ws.Range("D49")
$cD.Value2 = "'0.08311"
$cD.ClearFormats()
$ws.Range("E49").Value = '  +1.04%  '

$cD = $ws.Range("D50")
$cD.Value2 = "'134.46"
$cD.ClearFormats()
$ws.Range("E50").Value = '  +0.68%  '

$cD = $ws.Range("D51")
$cD.Value2 = "'1.224"
$cD.ClearFormats()
$ws.Range("E51").Value = '  +1.05%  '

Write-Host "Updated cryptos list"